# feat: db reset and equipe login
# Update progress tracking values for rows 7-11 (Temps passe / Reste a faire)
# and move the active selection to D24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E = Temps passe (en %), F = Reste a faire (en %), G = Avancement (en %) = (100 - F)
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 0

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 0

$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 0

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 0

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 0

# Move selection/active cell to D24 as in the committed workbook.
$ws.Range("D24").Select()
